$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 12 (date 44498)
$ws.Cells.Item(2, 4).Value = 44498
$ws.Cells.Item(2, 12).Value = "Segunda"
$ws.Cells.Item(2, 13).Value = 300
$ws.Cells.Item(2, 14).Value = 19000
$ws.Cells.Item(2, 15).Value = 20000
$ws.Cells.Item(2, 16).Value = 19500
$ws.Cells.Item(2, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(2, 19).Value = 2438
$ws.Cells.Item(2, 20).Value = 8

# Row 3 <- old row 11 (date 44162)
$ws.Cells.Item(3, 4).Value = 44162
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 2000
$ws.Cells.Item(3, 15).Value = 2100
$ws.Cells.Item(3, 16).Value = 2050
$ws.Cells.Item(3, 17).Value = "$/kilo (en caja de 14 kilos)"
$ws.Cells.Item(3, 19).Value = 2050
$ws.Cells.Item(3, 20).Value = 1

# Row 4 <- old row 10 (date 44516)
$ws.Cells.Item(4, 4).Value = 44516
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 19000
$ws.Cells.Item(4, 16).Value = 18500
$ws.Cells.Item(4, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(4, 19).Value = 2312
$ws.Cells.Item(4, 20).Value = 8

# Row 5 <- old row 2 (date 44488)
$ws.Cells.Item(5, 4).Value = 44488
$ws.Cells.Item(5, 12).Value = "Segunda"
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 17000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 17500
$ws.Cells.Item(5, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(5, 19).Value = 2188
$ws.Cells.Item(5, 20).Value = 8

# Row 6 <- old row 7 (date 44526)
$ws.Cells.Item(6, 4).Value = 44526
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 21000
$ws.Cells.Item(6, 15).Value = 21000
$ws.Cells.Item(6, 16).Value = 21000
$ws.Cells.Item(6, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(6, 19).Value = 2625
$ws.Cells.Item(6, 20).Value = 8

# Row 7 <- old row 15 (date 44523)
$ws.Cells.Item(7, 4).Value = 44523
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 400
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(7, 19).Value = 2688
$ws.Cells.Item(7, 20).Value = 8

# Row 8 <- old row 16 (date 44523)
$ws.Cells.Item(8, 4).Value = 44523
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 100
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(8, 19).Value = 2250
$ws.Cells.Item(8, 20).Value = 8

# Row 9 <- old row 18 (date 44495)
$ws.Cells.Item(9, 4).Value = 44495
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 270
$ws.Cells.Item(9, 14).Value = 19000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 19556
$ws.Cells.Item(9, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(9, 19).Value = 2444
$ws.Cells.Item(9, 20).Value = 8

# Row 10 <- old row 3 (date 44530)
$ws.Cells.Item(10, 4).Value = 44530
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 200
$ws.Cells.Item(10, 14).Value = 19000
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 19500
$ws.Cells.Item(10, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(10, 19).Value = 2438
$ws.Cells.Item(10, 20).Value = 8

# Row 11 <- old row 4 (date 44530)
$ws.Cells.Item(11, 4).Value = 44530
$ws.Cells.Item(11, 12).Value = "Segunda"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 16000
$ws.Cells.Item(11, 15).Value = 16000
$ws.Cells.Item(11, 16).Value = 16000
$ws.Cells.Item(11, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(11, 19).Value = 2000
$ws.Cells.Item(11, 20).Value = 8

# Row 12 <- old row 5 (date 44159)
$ws.Cells.Item(12, 4).Value = 44159
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 300
$ws.Cells.Item(12, 14).Value = 2000
$ws.Cells.Item(12, 15).Value = 2100
$ws.Cells.Item(12, 16).Value = 2050
$ws.Cells.Item(12, 17).Value = "$/kilo (en caja de 14 kilos)"
$ws.Cells.Item(12, 19).Value = 2050
$ws.Cells.Item(12, 20).Value = 1

# Row 14 <- old row 8 (date 44519)
$ws.Cells.Item(14, 4).Value = 44519
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 21000
$ws.Cells.Item(14, 15).Value = 22000
$ws.Cells.Item(14, 16).Value = 21500
$ws.Cells.Item(14, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(14, 19).Value = 2688
$ws.Cells.Item(14, 20).Value = 8

# Row 15 <- old row 9 (date 44519)
$ws.Cells.Item(15, 4).Value = 44519
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 18000
$ws.Cells.Item(15, 16).Value = 18000
$ws.Cells.Item(15, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(15, 19).Value = 2250
$ws.Cells.Item(15, 20).Value = 8

# Row 16 <- old row 6 (date 44512)
$ws.Cells.Item(16, 4).Value = 44512
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 19000
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 19500
$ws.Cells.Item(16, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(16, 19).Value = 2438
$ws.Cells.Item(16, 20).Value = 8

# Row 17 <- old row 14 (date 44509)
$ws.Cells.Item(17, 4).Value = 44509
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 19000
$ws.Cells.Item(17, 15).Value = 20000
$ws.Cells.Item(17, 16).Value = 19500
$ws.Cells.Item(17, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(17, 19).Value = 2438
$ws.Cells.Item(17, 20).Value = 8

# Row 18 <- old row 17 (date 44491)
$ws.Cells.Item(18, 4).Value = 44491
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 18000
$ws.Cells.Item(18, 15).Value = 19000
$ws.Cells.Item(18, 16).Value = 18500
$ws.Cells.Item(18, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(18, 19).Value = 2312
$ws.Cells.Item(18, 20).Value = 8
